# Laboratoriesystemer workbook update:
#   - The worksheet tab (and the matching defined name that points at it)
#     is dated; bump the date from 02-12-2025 to 05-12-2025.
# Renaming the sheet via the Name property also re-points any defined
# name / formula that references the sheet by its old name, exactly like
# a normal rename from the Excel UI does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Opdateret d. 05-12-2025"
